# Puliendo la estructura de los data frames de CPU y MEM
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab from procesosMEM to procesosCPU
$ws.Name = "procesosCPU"

# New data for columns B (CANTIDAD), C (MAX), D (MIN) for data rows 2-13
$data = @(
    @{ Row = 2;  B = 29; C = "0.00%";  D = "0.00%" },
    @{ Row = 3;  B = 29; C = "0.26%";  D = "0.06%" },
    @{ Row = 4;  B = 29; C = "0.04%";  D = "0.00%" },
    @{ Row = 5;  B = 29; C = "0.04%";  D = "0.00%" },
    @{ Row = 6;  B = 29; C = "0.04%";  D = "0.00%" },
    @{ Row = 7;  B = 29; C = "0.05%";  D = "0.00%" },
    @{ Row = 8;  B = 29; C = "0.06%";  D = "0.00%" },
    @{ Row = 9;  B = 29; C = "0.09%";  D = "0.00%" },
    @{ Row = 10; B = 29; C = "35.40%"; D = "0.22%" },
    @{ Row = 11; B = 29; C = "0.30%";  D = "0.22%" },
    @{ Row = 12; B = 29; C = "1.06%";  D = "0.67%" },
    @{ Row = 13; B = 29; C = "3.42%";  D = "0.69%" }
)

foreach ($item in $data) {
    $r = $item.Row

    # Column B is numeric (CANTIDAD)
    $ws.Cells.Item($r, 2).Value = $item.B

    # Columns C and D hold plain text that merely looks like percentages
    # (e.g. "0.07%"), so force text entry and then drop back to the
    # workbook's default "Normal" style to avoid leaving a custom number
    # format applied to the cell.
    $cellC = $ws.Cells.Item($r, 3)
    $cellC.NumberFormat = "@"
    $cellC.Value = $item.C
    $cellC.Style = "Normal"

    $cellD = $ws.Cells.Item($r, 4)
    $cellD.NumberFormat = "@"
    $cellD.Value = $item.D
    $cellD.Style = "Normal"
}
